$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 1, 31.61296591696135)
    3 = @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 0, 4.837881874639075)
    4 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759)
    5 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 15.68806981981553)
    6 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 1, 29.84159230404497)
    7 = @(0.3048080303191223, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 1, 1.001517020209437)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # B is column 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
